$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize existing "Creation date" text (col M, rows 2-183) from 4-digit
# to 2-digit year, e.g. "8/3/2018 19:27:00" -> "8/3/18 19:27:00".
$creationDates = @(
    "8/3/18 19:27:00",
    "8/3/18 19:28:00",
    "8/3/18 19:28:00",
    "8/3/18 19:28:00",
    "8/3/18 19:29:00",
    "8/3/18 19:29:00",
    "8/3/18 19:29:00",
    "8/3/18 19:31:00",
    "8/3/18 19:32:00",
    "8/3/18 19:32:00",
    "8/3/18 19:32:00",
    "8/3/18 19:32:00",
    "8/3/18 19:33:00",
    "8/3/18 19:33:00",
    "8/3/18 19:33:00",
    "8/3/18 19:34:00",
    "8/3/18 19:34:00",
    "8/3/18 19:35:00",
    "8/3/18 19:35:00",
    "8/3/18 19:35:00",
    "8/3/18 19:35:00",
    "8/3/18 19:36:00",
    "8/3/18 19:36:00",
    "8/3/18 19:36:00",
    "8/3/18 19:37:00",
    "8/3/18 19:37:00",
    "8/3/18 19:38:00",
    "8/3/18 19:38:00",
    "8/3/18 19:39:00",
    "8/3/18 19:41:00",
    "8/3/18 19:41:00",
    "8/3/18 19:41:00",
    "8/3/18 19:42:00",
    "8/3/18 19:42:00",
    "8/3/18 19:43:00",
    "8/3/18 19:43:00",
    "8/3/18 19:43:00",
    "8/3/18 19:43:00",
    "8/3/18 19:43:00",
    "8/3/18 19:43:00",
    "8/3/18 19:44:00",
    "8/3/18 19:44:00",
    "8/3/18 19:45:00",
    "8/3/18 19:45:00",
    "8/3/18 19:48:00",
    "8/3/18 19:49:00",
    "8/3/18 19:49:00",
    "8/3/18 19:49:00",
    "8/3/18 19:49:00",
    "8/3/18 19:49:00",
    "8/3/18 19:50:00",
    "8/3/18 19:50:00",
    "8/3/18 19:50:00",
    "8/3/18 19:50:00",
    "8/3/18 19:50:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:51:00",
    "8/3/18 19:52:00",
    "8/3/18 19:52:00",
    "8/3/18 19:52:00",
    "8/3/18 19:53:00",
    "8/3/18 19:53:00",
    "8/3/18 19:53:00",
    "8/3/18 19:55:00",
    "8/3/18 19:55:00",
    "8/3/18 19:56:00",
    "8/3/18 19:56:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:57:00",
    "8/3/18 19:58:00",
    "8/3/18 19:59:00",
    "8/3/18 19:59:00",
    "8/3/18 20:00:00",
    "8/3/18 20:01:00",
    "8/3/18 20:01:00",
    "8/3/18 20:01:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:02:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:03:00",
    "8/3/18 20:04:00",
    "8/3/18 20:04:00",
    "8/3/18 20:04:00",
    "8/3/18 20:05:00",
    "8/3/18 20:05:00",
    "8/3/18 20:05:00",
    "8/3/18 20:05:00",
    "8/3/18 20:05:00",
    "8/3/18 20:05:00",
    "8/3/18 20:06:00",
    "8/3/18 20:07:00",
    "8/3/18 20:08:00",
    "8/3/18 20:08:00",
    "8/3/18 20:08:00",
    "8/3/18 20:08:00",
    "8/3/18 20:10:00",
    "8/3/18 20:10:00",
    "8/3/18 20:10:00",
    "8/3/18 20:10:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:11:00",
    "8/3/18 20:12:00",
    "8/3/18 20:12:00",
    "8/3/18 20:12:00",
    "8/3/18 20:12:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:13:00",
    "8/3/18 20:14:00",
    "8/3/18 20:15:00",
    "8/3/18 20:32:00",
    "8/3/18 20:32:00",
    "8/3/18 20:32:00",
    "8/3/18 20:32:00",
    "8/3/18 20:32:00",
    "8/3/18 20:32:00",
    "10/15/18 14:43:00",
    "10/15/18 14:43:00",
    "10/15/18 14:43:00",
    "10/15/18 14:44:00",
    "10/15/18 14:44:00",
    "10/15/18 14:44:00",
    "10/15/18 14:46:00",
    "10/15/18 14:47:00",
    "10/15/18 14:49:00",
    "10/15/18 14:49:00",
    "10/15/18 14:52:00",
    "10/15/18 14:52:00",
    "10/15/18 14:52:00",
    "10/15/18 14:52:00",
    "11/8/18 11:38:00",
    "11/8/18 11:38:00",
    "11/8/18 11:38:00",
    "11/8/18 11:38:00",
    "11/8/18 11:39:00",
    "11/8/18 11:39:00",
    "11/8/18 11:39:00",
    "11/8/18 11:39:00",
    "11/8/18 11:39:00",
    "11/8/18 11:39:00",
    "11/8/18 11:40:00",
    "11/8/18 11:40:00",
    "11/8/18 11:40:00",
    "11/8/18 14:29:00",
    "11/12/18 13:33:00",
    "11/12/18 13:33:00",
    "11/12/18 13:33:00",
    "1/29/19 16:47:09"
)
for ($i = 0; $i -lt $creationDates.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $creationDates[$i]
}

# Append three new coded-segment rows (184-186), mirroring the row/style
# pattern already used throughout the sheet (row 3 as the style template,
# since its "Document name" column also equals "22121").
$ws.Range("A3:M3").Copy($ws.Range("A184:M184"))
$ws.Range("A3:M3").Copy($ws.Range("A185:M185"))
$ws.Range("A3:M3").Copy($ws.Range("A186:M186"))

# Row 184: Location:City -> Taipei
$ws.Range("E184").Value = "Location:City"
$ws.Range("F184").Value = "1: 1314"
$ws.Range("G184").Value = "1: 1319"
$ws.Range("H184").Value = 0
$ws.Range("I184").Value = "Taipei"
$ws.Range("J184").Value = 6
$ws.Range("K184").Value = 0.041531113726033089
$ws.Range("L184").Value = "emmamendelsohn"
$ws.Range("M184").Value = "8/22/19 14:19:16"

# Row 185: Location:Country -> Taiwan
$ws.Range("E185").Value = "Location:Country"
$ws.Range("F185").Value = "1: 1326"
$ws.Range("G185").Value = "1: 1331"
$ws.Range("H185").Value = 0
$ws.Range("I185").Value = "Taiwan"
$ws.Range("J185").Value = 6
$ws.Range("K185").Value = 0.041531113726033089
$ws.Range("L185").Value = "emmamendelsohn"
$ws.Range("M185").Value = "8/22/19 14:19:20"

# Row 186: Location:Hospital name -> Tri-Service General Hospital, National Defense Medical Center
$ws.Range("E186").Value = "Location:Hospital name"
$ws.Range("F186").Value = "1: 1213"
$ws.Range("G186").Value = "1: 1275"
$ws.Range("H186").Value = 0
$ws.Range("I186").Value = "Tri-Service General Hospital,  `nNational Defense Medical Center"
$ws.Range("J186").Value = 62
$ws.Range("K186").Value = 0.4291548418356752
$ws.Range("L186").Value = "emmamendelsohn"
$ws.Range("M186").Value = "8/22/19 14:19:48"

# Row heights: single-line rows render at 16pt, the wrapped two-line
# hospital-name row renders at 30pt (matches the sizing already used for
# other multi-line "Segment" cells in the sheet).
$ws.Rows.Item(184).RowHeight = 16
$ws.Rows.Item(185).RowHeight = 16
$ws.Rows.Item(186).RowHeight = 30
